# Update "Forecast Comparison" sheet with correct forecast output:
#  - insert a new "Week_Start_Date" column after "Week" (shifts ASIN.. right by one)
#  - shorten week labels from W01.."W16" to W1.."W16" (drop the leading zero)
#  - populate the new Week_Start_Date column with each week's start date
#  - store is_holiday_week as a proper boolean value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B ("Week_Start_Date"); everything from the old column B
# (ASIN) onward shifts one column to the right (B->C, C->D, ... I->J).
$ws.Columns.Item(2).Insert()
$ws.Range("B1").Value = "Week_Start_Date"

# Keep the dates as plain text (YYYY-MM-DD), not Excel date serials.
$ws.Range("B2:B17").NumberFormat = "@"

$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = "W" + ($i + 1)
    $ws.Range("B$row").Value = $weekStartDates[$i]
}

# is_holiday_week (now column J after the insert) should hold real booleans.
for ($row = 2; $row -le 17; $row++) {
    $cell = $ws.Cells.Item($row, 10)
    if ($cell.Value -eq 1) {
        $cell.Value = $true
    } else {
        $cell.Value = $false
    }
}

$ws.Range("A1").Select() | Out-Null
